$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.2
$ws.Range("H2").Value = 2.8
$ws.Range("I2").Value = 4.2
$ws.Range("J2").Value = 3.1
$ws.Range("L2").Value = 5
$ws.Range("S2").Value = 1.78
$ws.Range("T2").Value = 2.03
$ws.Range("Z2").Value = 21
$ws.Range("AG2").Value = 7
$ws.Range("AK2").Value = 41
$ws.Range("AO2").Value = 15
$ws.Range("AS2").Value = 501
$ws.Range("G3").Value = 2.88
$ws.Range("I3").Value = 3.25
$ws.Range("J3").Value = 3.75
$ws.Range("M3").Value = 1.2
$ws.Range("N3").Value = 4.33
$ws.Range("X3").Value = 11
$ws.Range("Z3").Value = 29
$ws.Range("AC3").Value = 4.33
$ws.Range("AQ3").Value = 67
$ws.Range("AU3").Value = 11
$ws.Range("AV3").Value = 101
$ws.Range("AW3").Value = 4.75
$ws.Range("AX3").Value = 23
$ws.Range("G4").Value = 2.3
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3.4
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.6
$ws.Range("X4").Value = 10
$ws.Range("AG4").Value = 8.5
$ws.Range("AH4").Value = 15
$ws.Range("AK4").Value = 29
$ws.Range("AX4").Value = 19
$ws.Range("BB4").Value = 251
$ws.Range("S5").Value = 1.33
$ws.Range("J6").Value = 1.83
$ws.Range("L6").Value = 7.5
$ws.Range("Q6").Value = 1.65
$ws.Range("R6").Value = 2.2
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 1.75
$ws.Range("W6").Value = 7.5
$ws.Range("Y6").Value = 9
$ws.Range("AA6").Value = 11
$ws.Range("AB6").Value = 26
$ws.Range("AC6").Value = 15
$ws.Range("AE6").Value = 21
$ws.Range("AI6").Value = 21
$ws.Range("AJ6").Value = 81
$ws.Range("AM6").Value = 401
$ws.Range("AU6").Value = 9
$ws.Range("BB6").Value = 301
$ws.Range("I7").Value = 4.1
$ws.Range("M7").Value = 1.14
$ws.Range("N7").Value = 5.5
$ws.Range("S7").Value = 1.67
$ws.Range("T7").Value = 2.1
$ws.Range("AA7").Value = 23
$ws.Range("AG7").Value = 7.5
$ws.Range("AJ7").Value = 51
$ws.Range("AT7").Value = 2.1
$ws.Range("G8").Value = 1.8
$ws.Range("I8").Value = 4.75
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("X8").Value = 7.5
$ws.Range("AD8").Value = 6.5
$ws.Range("AN8").Value = 3.6
$ws.Range("AO8").Value = 10
$ws.Range("AV8").Value = 67
$ws.Range("AZ8").Value = 101
$ws.Range("O9").Value = 1.44
$ws.Range("P9").Value = 2.63
$ws.Range("S10").Value = 1.25
$ws.Range("Q11").Value = 1.73
$ws.Range("R11").Value = 2.08
$ws.Range("S11").Value = 1.3
